$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.311.99"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "1.901.11"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.16%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.695"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +9.51%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "245.62"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("E7").Value = "  -0.17%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "40.80"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -4.17%  "
$ws.Range("E9").Value = "  +2.74%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "53.17"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +11.40%  "
$ws.Range("E11").Value = "  +2.05%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0995"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "2.176.95"
$ws.Range("E13").Value = "  -0.07%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "12.39"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("E15").Value = "  +2.24%  "
$ws.Range("D16").Value = "1.904.58"
$ws.Range("E16").Value = "  +0.29%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "4.81"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "35.312.24"
$ws.Range("E18").Value = "  -0.92%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "72.24"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  +0.98%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "240.78"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.55%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "12.61"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.15%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.79"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("E24").Value = "  -0.10%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.30"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("E26").Value = "  +7.72%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "168.23"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.76%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "8.60"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.15%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "18.34"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.95%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.131"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +3.84%  "
$ws.Range("D31").Value = "4.141.91"
$ws.Range("E31").Value = "  +21.32%  "
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("E34").Value = "  -0.11%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.922"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -3.78%  "
$ws.Range("E36").Value = "  -0.09%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.80"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.24%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.48"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +10.71%  "
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("E40").Value = "  +10.89%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0210"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("E42").Value = "  -1.03%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "16.00"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +4.00%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "89.71"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("D45").Value = "1.350.37"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("E46").Value = "  +2.85%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "12.65"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.50%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.43"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  +0.64%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "46.06"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.35%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "6.51"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.08%  "
